$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("D6").Value = "2016-03-09 14:09:58"
$wsDe.Range("D6").Value = "2016-03-09 14:10:01"
